$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    # Find the shape holding the table with the largest row count on this
    # slide (there's normally just one table per slide, but this keeps the
    # logic robust if a slide ever has more than one).
    $maxRowCount = 0
    $maxRowIndex = -1
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shape = $s.Shapes.Item($j)
        if ($shape.HasTable) {
            $rowCount = $shape.Table.Rows.Count
            if ($rowCount -gt $maxRowCount) {
                $maxRowCount = $rowCount
                $maxRowIndex = $j
            }
        }
    }

    if ($maxRowIndex -gt 0) {
        $s.Shapes.Item($maxRowIndex).Delete()
    }
}
